# "fixes the help menu"
# The listing sheet had a stray header row ("RPO" / "States") sitting above
# the actual RPO/States data table. Remove that header row so the data
# (and the query-table backed "ertac_rpo_listing" range) starts at row 1,
# then point the defined name at the new A1:B4 range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the header row (row 1: "RPO" / "States"); this shifts the
# remaining data rows up by one and drops the now-empty former row 5.
$ws.Rows.Item(1).Delete()

# Re-point the defined name used by the query table to the new data extent.
$wb.Names.Item("ertac_rpo_listing").RefersTo = "=Sheet1!`$A`$1:`$B`$4"

# Leave the active selection where the author left it.
$ws.Range("H13").Select() | Out-Null
